$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT (no auto number/date coercion),
# matching the original inlineStr cells, by temporarily forcing a text
# number format, then resetting the cell style back to the default/unstyled
# look (copied from an always-unstyled neighbor) so no stray style is left
# behind in the saved workbook.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $ws.Cells.Item(1, 1).Style
}

$ws.Cells.Item(2, 4).Value = '58.562.80'
$ws.Cells.Item(3, 4).Value = '2.568.40'
Set-TextValue 5 4 '509.88'
Set-TextValue 6 4 '146.40'
Set-TextValue 7 4 '0.998'
Set-TextValue 8 4 '0.573'
$ws.Cells.Item(9, 4).Value = '2.585.74'
Set-TextValue 10 4 '6.31'
Set-TextValue 11 4 '0.103'
Set-TextValue 12 4 '0.336'
$ws.Cells.Item(15, 4).Value = '58.520.16'
Set-TextValue 16 4 '21.21'
Set-TextValue 17 4 '0.0000137'
$ws.Cells.Item(18, 4).Value = '2.580.80'
Set-TextValue 19 4 '4.56'
Set-TextValue 20 4 '346.67'
Set-TextValue 21 4 '10.33'
Set-TextValue 24 4 '60.78'
Set-TextValue 25 4 '0.418'
Set-TextValue 27 4 '0.161'
$ws.Cells.Item(28, 4).Value = '2.687.00'
$ws.Cells.Item(29, 4).Value = '0.0₃0815'
Set-TextValue 30 4 '7.02'
Set-TextValue 32 4 '6.08'
Set-TextValue 33 4 '18.78'
Set-TextValue 34 4 '149.76'
Set-TextValue 36 4 '0.942'
Set-TextValue 37 4 '3.99'
Set-TextValue 39 4 '0.857'
Set-TextValue 40 4 '36.14'
Set-TextValue 41 4 '293.92'
Set-TextValue 42 4 '1.40'
Set-TextValue 43 4 '3.58'
Set-TextValue 44 4 '0.0996'
Set-TextValue 45 4 '0.997'
Set-TextValue 46 4 '0.610'
Set-TextValue 47 4 '0.0539'
Set-TextValue 48 4 '19.13'

$ws.Cells.Item(2, 5).Value = '  -4.23%  '
$ws.Cells.Item(3, 5).Value = '  -3.47%  '
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$ws.Cells.Item(5, 5).Value = '  -4.26%  '
$ws.Cells.Item(6, 5).Value = '  -6.59%  '
$ws.Cells.Item(7, 5).Value = '  +0.09%  '
$ws.Cells.Item(8, 5).Value = '  -3.09%  '
$ws.Cells.Item(9, 5).Value = '  -3.27%  '
$ws.Cells.Item(10, 5).Value = '  -4.50%  '
$ws.Cells.Item(11, 5).Value = '  -5.47%  '
$ws.Cells.Item(12, 5).Value = '  -5.10%  '
$ws.Cells.Item(13, 5).Value = '  -0.72%  '
$ws.Cells.Item(15, 5).Value = '  -4.27%  '
$ws.Cells.Item(16, 5).Value = '  -3.93%  '
$ws.Cells.Item(17, 5).Value = '  -4.51%  '
$ws.Cells.Item(18, 5).Value = '  -3.44%  '
$ws.Cells.Item(19, 5).Value = '  -4.89%  '
$ws.Cells.Item(20, 5).Value = '  -2.64%  '
$ws.Cells.Item(21, 5).Value = '  -3.64%  '
$ws.Cells.Item(22, 5).Value = '  -3.56%  '
$ws.Cells.Item(23, 5).Value = '  -0.02%  '
$ws.Cells.Item(24, 5).Value = '  -1.34%  '
$ws.Cells.Item(25, 5).Value = '  -3.39%  '
$ws.Cells.Item(26, 5).Value = '  -0.18%  '
$ws.Cells.Item(27, 5).Value = '  -4.65%  '
$ws.Cells.Item(28, 5).Value = '  -3.32%  '
$ws.Cells.Item(29, 5).Value = '  -5.85%  '
$ws.Cells.Item(30, 5).Value = '  -5.48%  '
$ws.Cells.Item(31, 5).Value = '  -0.04%  '
$ws.Cells.Item(32, 5).Value = '  -1.82%  '
$ws.Cells.Item(33, 5).Value = '  -4.16%  '
$ws.Cells.Item(34, 5).Value = '  -0.10%  '
$ws.Cells.Item(35, 5).Value = '  -5.30%  '
$ws.Cells.Item(36, 5).Value = '  +6.42%  '
$ws.Cells.Item(37, 5).Value = '  -3.98%  '
$ws.Cells.Item(38, 5).Value = '  -5.43%  '
$ws.Cells.Item(39, 5).Value = '  -6.55%  '
$ws.Cells.Item(40, 5).Value = '  -2.14%  '
$ws.Cells.Item(41, 5).Value = '  -4.95%  '
$ws.Cells.Item(42, 5).Value = '  -6.68%  '
$ws.Cells.Item(43, 5).Value = '  -6.34%  '
$ws.Cells.Item(44, 5).Value = '  -2.72%  '
$ws.Cells.Item(45, 5).Value = '  +0.01%  '
$ws.Cells.Item(46, 5).Value = '  -6.29%  '
$ws.Cells.Item(47, 5).Value = '  -4.76%  '
$ws.Cells.Item(48, 5).Value = '  -6.54%  '

# Rows 49-51 reordered/renamed coins
$ws.Cells.Item(49, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 49 4 '10.27'
$ws.Cells.Item(49, 5).Value = '  -0.79%  '

$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 50 4 '0.0229'
$ws.Cells.Item(50, 5).Value = '  -4.78%  '

$ws.Cells.Item(51, 2).Value = 'RenderToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 51 4 '4.64'
$ws.Cells.Item(51, 5).Value = '  -7.32%  '

